$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Google Drive direct-download URLs (export=download& removed) for column C (image),
# rows 2-25, corresponding to the 24 "image" shared-string entries that were edited.
$newImageUrls = @(
    "https://drive.google.com/uc?id=1HMYZcMib530mPnXfWylSLhhELM1fVUCV=FILE_ID",
    "https://drive.google.com/uc?id=1XpHtsMXASzE0c30eol_L4v5yxtVv5M6D",
    "https://drive.google.com/uc?id=1gvY2M2D5uQXOLPq8_nrTDcuhnOgDt1mv",
    "https://drive.google.com/uc?id=1RUoXaI0tCiTUVJH1CwDNBmrGHG3B38ce",
    "https://drive.google.com/uc?id=1XgVEvBeEw9mmH4Sg-U0lIynjI3X8uA8y",
    "https://drive.google.com/uc?id=1q90SEhpTxSGQDiY652adYF2Ej-1DVJ1u",
    "https://drive.google.com/uc?id=1we3nj6T0gaSAfeETAA2_ffkdbHKTZKe",
    "https://drive.google.com/uc?id=1N16AeaEf9B3jbAxah2h3q50F7TeMCJY3",
    "https://drive.google.com/uc?id=1hJhePKwM83WWBndemzHdwgPPgeim6aVH",
    "https://drive.google.com/uc?id=1cAUO7NQ3fsnIvx37Yc1HO5OMfD8UlLPs",
    "https://drive.google.com/uc?id=1W_Ol-b1SIJxAIrKhVFi9_gniq3DFH2Be",
    "https://drive.google.com/uc?id=1g3rRetLejsHMDw1AK6-BkPnbVMHyVaPc",
    "https://drive.google.com/uc?id=1rP2J6qkbwII2NNKLyf76KLStJxImiEjZ",
    "https://drive.google.com/uc?id=1mYl5_kMbOk07srSM5AC6Fd1D2a5t6Mvz",
    "https://drive.google.com/uc?id=1p_y9MXF2Lez29G1GCFHRCNLq34wGQsCP",
    "https://drive.google.com/uc?id=1PNshTU_Gt27rfePDfxLwEZIoa-UPN4i-",
    "https://drive.google.com/uc?id=1iuQNBlX5ULFJuj5Mh1Y8OHQTqJc4df4D",
    "https://drive.google.com/uc?id=1FMnHlklBjt5oe3Aul1YirbwhnrGFoTII",
    "https://drive.google.com/uc?id=1ihdSI889_Y1wxO2DArrwKX_OLV-McKsk",
    "https://drive.google.com/uc?id=16KolvaIcKm6Bds2EMEUL480WKA515TtF",
    "https://drive.google.com/uc?id=1oxs_MMNisWSe9WHpAzi2Gp0GS3Ia-0Lo",
    "https://drive.google.com/uc?id=1dN9gMbnypUO_-PNlOS73LcW61x5-AodN",
    "https://drive.google.com/uc?id=1yHht1DyjUjah-xEmcsu5N9ydwB1hKi-h",
    "https://drive.google.com/uc?id=1McTmDZSsTXGVksTI6rWVejaSv3xxYOnh"
)

for ($i = 0; $i -lt $newImageUrls.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newImageUrls[$i]
}
